$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) FCT 1A: change selection from A6 to the whole used range A1:N6
# ---------------------------------------------------------------------------
$ws1A = $wb.Worksheets.Item("FCT 1A ")
$ws1A.Activate()
$ws1A.Range("A1:N6").Select()

# ---------------------------------------------------------------------------
# 2) FCT 2B: add SiO2 (C,D) and Cl (I,J) measurement columns, then move the
#    selection to J3:J5 (this also removes the tabSelected flag once the
#    new "Dunk" sheet becomes active below)
# ---------------------------------------------------------------------------
$ws2B = $wb.Worksheets.Item("FCT 2B")
$ws2B.Activate()

$ws2B.Range("C2").Value = 42.6
$ws2B.Range("D2").Value = 56.2

$ws2B.Range("C3").Value = 41.6
$ws2B.Range("D3").Value = 56.4
$ws2B.Range("I3").Value = 35
$ws2B.Range("J3").Value = 60

$ws2B.Range("C4").Value = 43.8
$ws2B.Range("D4").Value = 53
$ws2B.Range("I4").Value = 37.4
$ws2B.Range("J4").Value = 61

$ws2B.Range("C5").Value = 46.4
$ws2B.Range("D5").Value = 49.4
$ws2B.Range("I5").Value = 39.6
$ws2B.Range("J5").Value = 62.6

$ws2B.Range("C6").Value = 48.2
$ws2B.Range("D6").Value = 49.8

$ws2B.Range("J3:J5").Select()

# ---------------------------------------------------------------------------
# 3) Add a new worksheet "Dunk" at the end of the workbook with the ICR data
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$dunk = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$dunk.Name = "Dunk"

# Row labels (A2:A5) first, so the new shared strings are created in the
# same order as in the target workbook
$dunk.Range("A2").Value = "D1"
$dunk.Range("A3").Value = "D2"
$dunk.Range("A4").Value = "D3"
$dunk.Range("A5").Value = "D4"

# Header row
$dunk.Range("A1").Value = "name"
$dunk.Range("B1").Value = "tid"
$dunk.Range("C1").Value = "FCT_1A SiO2"
$dunk.Range("D1").Value = "FCT_1B SiO2"
$dunk.Range("E1").Value = "FCT_2A SiO2"
$dunk.Range("F1").Value = "FCT_2B SiO2"
$dunk.Range("G1").Value = "FCT_1A Cl"
$dunk.Range("H1").Value = "FCT_1BCl"
$dunk.Range("I1").Value = "FCT_2A Cl"
$dunk.Range("J1").Value = "FCT_2B Cl"

# tid column formulas
$dunk.Range("B2").Formula = "=2100/1002"
$dunk.Range("B3").Formula = "=4200/1002"
$dunk.Range("B4").Formula = "=3600/1002"
$dunk.Range("B5").Formula = "=8800/1002"

# Data values
$dunk.Range("C2").Value = 41.4
$dunk.Range("D2").Value = 38.4
$dunk.Range("E2").Value = 37.2
$dunk.Range("F2").Value = 35.1

$dunk.Range("C3").Value = 43
$dunk.Range("D3").Value = 40.3
$dunk.Range("E3").Value = 39.5
$dunk.Range("F3").Value = 35.6
$dunk.Range("G3").Value = 87.4
$dunk.Range("H3").Value = 81.4
$dunk.Range("I3").Value = 71.4
$dunk.Range("J3").Value = 59.2

$dunk.Range("C4").Value = 47.2
$dunk.Range("D4").Value = 41.3
$dunk.Range("E4").Value = 41.3
$dunk.Range("F4").Value = 36.7

$dunk.Range("C5").Value = 49.6
$dunk.Range("D5").Value = 46.5
$dunk.Range("E5").Value = 46.3
$dunk.Range("F5").Value = 38.3
$dunk.Range("G5").Value = 72.2
$dunk.Range("H5").Value = 87.4
$dunk.Range("I5").Value = 61.2
$dunk.Range("J5").Value = 69.8

# Column C width (best-fit in the original file)
$dunk.Columns.Item(3).ColumnWidth = 11

# Final selection on the new sheet
$dunk.Range("L10").Select()
